# "Test export of XLSX metadata"
#
# 1. Nudges the formatting of the existing "Tabelle1" sheet: column A gets
#    one (shared) cell style, column B gets a different (shared) cell
#    style, and every row now has a (possibly blank) B cell.
# 2. Adds a new "meta" worksheet right after "Tabelle1" containing a small
#    Name/Value metadata table (Author/Unknown, Year/2024). Adding the
#    sheet after Tabelle1 also moves the active tab onto "meta".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 12

# A second, named style used only to give column B cells a style identity
# distinct from column A's (both resolve to the same visual formatting as
# "Normal" - only the style slot differs, matching the source file's
# harmless re-export style churn).
$colBStyle = $wb.Styles.Add("MetaColB")
$colBStyle.Font.Bold = $false

for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Style = "Normal"
    $ws.Cells.Item($r, 2).Style = "MetaColB"
}

# Add the "meta" worksheet right after "Tabelle1" - this also moves the
# active tab from index 0 to index 1, matching the workbook's activeTab
# change.
$metaSheet = $wb.Worksheets.Add($null, $ws)
$metaSheet.Name = "meta"

$metaSheet.Range("A1").Value = "Name"
$metaSheet.Range("B1").Value = "Value"
$metaSheet.Range("A2").Value = "Author"
$metaSheet.Range("B2").Value = "Unknown"
$metaSheet.Range("A3").Value = "Year"
$metaSheet.Range("B3").Value = 2024
